# This script updates the crypto price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are purely numeric-looking text (e.g. '5.20', '1.00') that must be
# preserved exactly (including trailing zeros) as text, not converted to numbers.
# Temporarily mark those cells as Text before assigning, then restore the default
# 'Normal' cell style so formatting matches the rest of the sheet.
$textCells = @('D5','D6','D8','D12','D14','D15','D19','D20','D21','D22','D23','D25','D26','D27','D30','D31','D32','D33','D34','D35','D37','D39','D40','D41','D42','D43','D45','D46','D47','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.331.87'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '2.450.99'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '577.93'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '143.42'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '2.448.64'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('D12').Value = '5.20'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').Value = '26.32'
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').Value = '0.0000175'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = '2.896.28'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '62.177.10'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '2.444.76'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = '10.88'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '327.85'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = '4.12'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = '1.96'
$ws.Range('E23').Value = '  -6.52%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '65.64'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '9.33'
$ws.Range('E26').Value = '  +3.73%  '
$ws.Range('D27').Value = '586.36'
$ws.Range('E27').Value = '  -5.66%  '
$ws.Range('D28').Value = '2.571.85'
$ws.Range('D29').Value = '0.0₃0956'
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').Value = '1.88'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '0.136'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').Value = '4.92'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').Value = '1.44'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').Value = '153.04'
$ws.Range('E39').Value = '  +3.99%  '
$ws.Range('D40').Value = '5.32'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').Value = '18.41'
$ws.Range('E41').Value = '  -2.11%  '
$ws.Range('D42').Value = '43.19'
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('D43').Value = '1.72'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  -4.13%  '
$ws.Range('D46').Value = '3.65'
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('D47').Value = '141.69'
$ws.Range('E47').Value = '  -3.06%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0250'
$ws.Range('E48').Value = '  +11.09%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.606'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('D50').Value = '0.0521'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').Value = '19.73'
$ws.Range('E51').Value = '  -4.35%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
